$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.510.94"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.418.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.420.62"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.007.84"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.533.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.414.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.13"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.61"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0760"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.848.16"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.89"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.43"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.93%  "
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.768"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.48"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.74%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.58"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.855"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.97%  "
